# Update Groupe, Noms, Thème et signature
$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1) Fill in the three empty "informations generales" table cells ---
$t = $d.Tables(1)

# Row 1 "Groupe" -> "OnlyUpSàrl" (flagged as a spell-check exception)
$cellGroupe = $t.Cell(1, 2)
$xmlGroupe = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>OnlyUpSàrl</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$cellGroupe.Range.InsertXML($xmlGroupe)

# Row 2 "Membres" -> "Eliott – Dioussé - Nicola"
$cellMembres = $t.Cell(2, 2)
$xmlMembres = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Eliott – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dioussé</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Nicola</w:t></w:r></w:p>
'@
$cellMembres.Range.InsertXML($xmlMembres)

# Row 3 "Thème choisi" -> " Immeuble pour habitations et bureaux"
$cellTheme = $t.Cell(3, 2)
$xmlTheme = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> Immeuble pour habitations et bureaux</w:t></w:r></w:p>
'@
$cellTheme.Range.InsertXML($xmlTheme)

# --- 2) Append the sign-off block at the end of the document body ---
$emptyCorpsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Corpsdetexte"/></w:pPr></w:p>'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($d.Paragraphs.Count)
$p1.Range.InsertXML($emptyCorpsXml)

$p1 = $d.Paragraphs($d.Paragraphs.Count)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($d.Paragraphs.Count)
$p2.Range.InsertXML($emptyCorpsXml)

$p2 = $d.Paragraphs($d.Paragraphs.Count)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($d.Paragraphs.Count)
$p3.Range.InsertXML($emptyCorpsXml)

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($d.Paragraphs.Count)
$xmlSignature = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Corpsdetexte"/></w:pPr><w:r><w:t xml:space="preserve">Lu et approuvé, </w:t></w:r><w:r><w:t>23.01.24, Scherrer Eliott</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$p4.Range.InsertXML($xmlSignature)
